{"js": "// Split the single bibliography paragraph (currently one run holding all\n// five numbered references concatenated together) into the same run with\n// the five references separated by manual line breaks (<w:br/>), i.e.\n// Word.BreakType.line, instead of being crammed one after another.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the bibliography paragraph robustly by its distinctive content\n// rather than by a hard-coded index.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"Peddy, S. The art of mentoring\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Bibliography paragraph not found\");\n}\n\n// The five references, in order, exactly as they appear in the original\n// run (split right before each \"[n]\" marker).\nconst refs = [\n  \"[1] Peddy, S. The art of mentoring \\u2013 Lead, follow and get out of the way. Houston: Bullion Books, 2001.\",\n  \"[2] Zachary, L. J. The Mentor\\u2019s Guide. San Francisco: Jossey-Bass Publishers, 2000. Pereira, A. Modelos de desenvolvimento do jovem adulto e promo\u00e7\u00e3o do bem-estar em estudantes do ensino superior. In: Programa de Monitoriza\u00e7\u00e3o e Tutorado: oito anos a promover a integra\u00e7\u00e3o e o sucesso acad\u00e9mico no IST. Lisboa: IST Press, 2011. p. 19-27.\",\n  \"[3] Mueller, S. Electronic mentoring as an example for the use of information and communications technology in engineering education. European Journal of Engineering Education, 2004.\",\n  \"[4] Kaul, S. Triangulated Mentorship of Engineering Students - Leveraging Peer Mentoring and Vertical Integration, Global Journal of Engineering Education, v. 21, p. 14-23,2019.\",\n  \"[5] Diretrizes Curriculares Nacionais para os cursos de gradua\u00e7\u00e3o em Engenharia. Minist\u00e9rio da Educa\u00e7\u00e3o. CNE/CES, 2019.\"\n];\n\n// Word represents a manual line break (\"\\v\", i.e. <w:br/>) with the\n// vertical-tab character \"\\u000b\" when reading/writing Range.text. Joining\n// the references with that character and replacing the whole paragraph's\n// text reproduces \"<w:t/><w:br/>\" pairs inside the single run.\nconst newText = refs.join(\"\\u000b\");\n\ntarget.insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Split the single bibliography paragraph (currently one run holding all\n# five numbered references concatenated together) so the five references\n# are separated by manual line breaks (vertical-tab / Chr(11), which Word\n# round-trips to <w:br/>) instead of being crammed one after another.\n\n$d = $word.ActiveDocument\n\n# Locate the bibliography paragraph robustly by its distinctive content\n# rather than by a hard-coded index.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Peddy, S. The art of mentoring*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw \"Bibliography paragraph not found\"\n}\n\n# The five references, in order, exactly as they appear in the original\n# run (split right before each \"[n]\" marker).\n$refs = @(\n  \"[1] Peddy, S. The art of mentoring \u2013 Lead, follow and get out of the way. Houston: Bullion Books, 2001.\",\n  \"[2] Zachary, L. J. The Mentor\u2019s Guide. San Francisco: Jossey-Bass Publishers, 2000. Pereira, A. Modelos de desenvolvimento do jovem adulto e promo\u00e7\u00e3o do bem-estar em estudantes do ensino superior. In: Programa de Monitoriza\u00e7\u00e3o e Tutorado: oito anos a promover a integra\u00e7\u00e3o e o sucesso acad\u00e9mico no IST. Lisboa: IST Press, 2011. p. 19-27.\",\n  \"[3] Mueller, S. Electronic mentoring as an example for the use of information and communications technology in engineering education. European Journal of Engineering Education, 2004.\",\n  \"[4] Kaul, S. Triangulated Mentorship of Engineering Students - Leveraging Peer Mentoring and Vertical Integration, Global Journal of Engineering Education, v. 21, p. 14-23,2019.\",\n  \"[5] Diretrizes Curriculares Nacionais para os cursos de gradua\u00e7\u00e3o em Engenharia. Minist\u00e9rio da Educa\u00e7\u00e3o. CNE/CES, 2019.\"\n)\n\n# Word represents a manual line break (\"\\v\", i.e. <w:br/>) with Chr(11)\n# (vertical tab) inside Range.Text. Joining the references with that\n# character and writing it back to the paragraph's range (minus its\n# trailing paragraph mark) reproduces \"<w:t/><w:br/>\" pairs inside the\n# single run.\n$r = $target.Range\n$r.End = $r.End - 1\n$r.Text = $refs -join [char]0x000b\n"}
